$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set for the arts listing (rows 2..59 after refresh; old row 60 removed)
$data = @(
    @(3, 'OPERATIONAL', 'Art On Main Gallery And Gifts', 5, 11),
    @(54, 'OPERATIONAL', 'Art To Remember', 3.9, 7),
    @(14, 'OPERATIONAL', 'ArtMix', 5, 8),
    @(0, 'OPERATIONAL', 'Arts Council of Indianapolis', 4.8, 12),
    @(7, 'OPERATIONAL', 'Arts for Lawrence', 4.9, 19),
    @(13, 'OPERATIONAL', 'Arts for Learning', 0, 0),
    @(35, 'OPERATIONAL', 'Carmel International Arts Festival', 4.8, 18),
    @(24, 'OPERATIONAL', 'Cat Head Press: Printshop and Artist Cooperative', 5, 26),
    @(10, 'OPERATIONAL', 'Christel De Haan Fine Arts Center', 4.6, 15),
    @(30, 'OPERATIONAL', 'Clowes Memorial Hall', 4.7, 797),
    @(51, 'OPERATIONAL', 'Creative Energy Arts', 0, 0),
    @(17, 'OPERATIONAL', 'Creative Expressions Arts', 0, 0),
    @(45, 'OPERATIONAL', 'Cultural Arts Gallery', 5, 1),
    @(31, 'OPERATIONAL', 'Dance Arts', 4.8, 10),
    @(55, 'OPERATIONAL', 'Dance Magic Performing Arts', 0, 0),
    @(44, 'OPERATIONAL', 'Edison School Of The Arts', 3.6, 36),
    @(11, 'OPERATIONAL', 'Gallery 924', 4.9, 15),
    @(8, 'OPERATIONAL', 'Garfield Park Arts Center (GPAC)', 4.6, 119),
    @(27, 'OPERATIONAL', 'Harrison Center', 4.8, 134),
    @(38, 'OPERATIONAL', 'Herron School of Art and Design', 4.8, 20),
    @(33, 'OPERATIONAL', 'Herron School of Art and Design - Eskenazi Fine Arts Center', 4.5, 2),
    @(15, 'OPERATIONAL', 'High Frequency Arts', 5, 3),
    @(46, 'OPERATIONAL', 'Ignition Arts, LLC', 4.9, 8),
    @(16, 'OPERATIONAL', 'Incredible Art Department', 0, 0),
    @(2, 'OPERATIONAL', 'Indiana Arts Commission', 0, 0),
    @(37, 'OPERATIONAL', 'Indiana Fine Arts Academy', 4, 4),
    @(48, 'OPERATIONAL', 'Indiana Performing Arts Centre', 4.8, 6),
    @(4, 'OPERATIONAL', 'Indianapolis Art Center', 4.7, 112),
    @(57, 'OPERATIONAL', 'Indianapolis Arts Chorale', 0, 0),
    @(19, 'OPERATIONAL', 'Indianapolis Arts Services', 5, 1),
    @(18, 'OPERATIONAL', 'Indianapolis Artsgarden', 4.6, 146),
    @(52, 'OPERATIONAL', 'Jazz Arts Society of In Inc', 4, 1),
    @(20, 'OPERATIONAL', 'Korka International Arts', 0, 0),
    @(43, 'OPERATIONAL', 'Landmark Keystone Art Cinema', 4.4, 513),
    @(12, 'OPERATIONAL', 'Magdalena Gallery of Arts', 5, 1),
    @(32, 'OPERATIONAL', 'Mars Hill Arts Center', 4.9, 15),
    @(26, 'OPERATIONAL', 'Newfields', 4.7, 3896),
    @(9, 'OPERATIONAL', 'Nickel Plate Arts', 4.9, 41),
    @(34, 'OPERATIONAL', 'Palette Art Studio', 5, 7),
    @(58, 'OPERATIONAL', 'Peddler''s Marketplace Folk Art', 0, 0),
    @(47, 'OPERATIONAL', 'Performing Arts', 0, 0),
    @(50, 'OPERATIONAL', 'Pike Performing Arts Center', 4.6, 180),
    @(56, 'OPERATIONAL', 'Red Barn Arts Collective', 5, 1),
    @(22, 'OPERATIONAL', 'Ripple Mobile Arts', 5, 7),
    @(23, 'OPERATIONAL', 'Schrott Center for the Arts', 4.8, 106),
    @(1, 'OPERATIONAL', 'Sho Arts', 5, 1),
    @(42, 'OPERATIONAL', 'Stallion Art Decor Inc', 4.1, 8),
    @(29, 'OPERATIONAL', 'Studio Alchemy LLC', 5, 5),
    @(25, 'OPERATIONAL', 'The Art Studio of Carmel', 0, 0),
    @(6, 'OPERATIONAL', 'The Carmel Arts Council Children’s Art Gallery', 5, 1),
    @(28, 'OPERATIONAL', 'The Center for the Performing Arts', 4.7, 1193),
    @(5, 'OPERATIONAL', 'The Murphy Art Center', 5, 1),
    @(53, 'OPERATIONAL', 'The Palladium at the Center for the Performing Arts', 4.9, 197),
    @(36, 'OPERATIONAL', 'United Art & Education', 4.5, 132),
    @(39, 'OPERATIONAL', 'Warren Performing Arts Center', 4.5, 374),
    @(49, 'OPERATIONAL', 'Zionsville Performing Arts Center', 4.7, 175),
    @(40, 'OPERATIONAL', 'carmel academy of the arts', 0, 0),
    @(59, 'OPERATIONAL', 'iDREAM Music And Arts Center', 0, 0)
)

$rowCount = $data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Remove the now-obsolete 60th row (dataset shrank by one row)
$ws.Rows.Item(60).Delete()
